$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "RET-36547"
$ws.Range("B2").Value = "Bismillah Mobile Shop 2"
$ws.Range("C2").Value = "Jewel Telecom"
$ws.Range("D2").Value = "N/A"
$ws.Range("E2").Value = "N/A"
$ws.Range("F2").Value = "N/A"

$ws.Range("D13").Select()
